$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-coerced to a number by
# Excel's usual "typed into a cell" parsing get pre-formatted as Text so the
# literal string (matching the original inlineStr cell) is preserved.
$textCells = @(
    "D5",
    "D7",
    "D10",
    "D12",
    "D13",
    "D15",
    "D16",
    "D19",
    "D22",
    "D24",
    "D27",
    "D28",
    "D29",
    "D33",
    "D37",
    "D38",
    "D39",
    "D45",
    "D47",
    "D48",
    "D51",
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '37.181.78'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '2.001.86'
$ws.Range("E3").Value = '  +2.08%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '246.08'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("E6").Value = '  +1.84%  '
$ws.Range("D7").Value = '59.78'
$ws.Range("E7").Value = '  +1.58%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +2.42%  '
$ws.Range("D10").Value = '0.0806'
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("D12").Value = '15.04'
$ws.Range("E12").Value = '  +5.72%  '
$ws.Range("D13").Value = '22.44'
$ws.Range("E13").Value = '  +5.75%  '
$ws.Range("D14").Value = '2.296.73'
$ws.Range("E14").Value = '  +2.18%  '
$ws.Range("D15").Value = '0.845'
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = '5.43'
$ws.Range("E16").Value = '  +2.51%  '
$ws.Range("D17").Value = '2.004.74'
$ws.Range("E17").Value = '  +2.27%  '
$ws.Range("D18").Value = '37.112.05'
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("D19").Value = '70.32'
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").Value = '0.0₃0863'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").Value = '230.35'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = '2.46'
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("E26").Value = '  +2.92%  '
$ws.Range("D27").Value = '0.141'
$ws.Range("E27").Value = '  +2.20%  '
$ws.Range("D28").Value = '163.51'
$ws.Range("E28").Value = '  +1.66%  '
$ws.Range("D29").Value = '19.65'
$ws.Range("E29").Value = '  +1.02%  '
$ws.Range("E30").Value = '  +11.18%  '
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("E32").Value = '  +1.29%  '
$ws.Range("D33").Value = '0.0655'
$ws.Range("E33").Value = '  +6.80%  '
$ws.Range("E34").Value = '  +2.00%  '
$ws.Range("E35").Value = '  +4.78%  '
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("D37").Value = '1.81'
$ws.Range("E37").Value = '  +2.09%  '
$ws.Range("D38").Value = '3.30'
$ws.Range("E38").Value = '  -6.35%  '
$ws.Range("D39").Value = '5.39'
$ws.Range("E39").Value = '  -1.06%  '
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("E41").Value = '  +0.83%  '
$ws.Range("E42").Value = '  +2.05%  '
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("E44").Value = '  +5.41%  '
$ws.Range("D45").Value = '90.95'
$ws.Range("E45").Value = '  +3.21%  '
$ws.Range("D46").Value = '1.373.80'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = '1.04'
$ws.Range("E47").Value = '  +1.83%  '
$ws.Range("D48").Value = '7.44'
$ws.Range("E48").Value = '  +4.17%  '
$ws.Range("E49").Value = '  +12.54%  '
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("D51").Value = '46.17'
$ws.Range("E51").Value = '  +4.78%  '
